# Legs and Room 3 Sesi 2 Update!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# L3 row (row 6): D (beta 45) value update, dependent H/K formulas auto-recalc
$ws.Range("D6").Value = 1650

# R1 row (row 7): D (beta 45) value update, dependent H/K formulas auto-recalc
$ws.Range("D7").Value = 1350

# R2 row (row 8): E (alpha 0) value update, dependent M formula auto-recalc
$ws.Range("E8").Value = 1870

# Update the active selection to match the saved cursor position
$ws.Range("F14").Select() | Out-Null
